$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "line": add two new double lines (Line6, Line7) as rows 8 and 9
# ---------------------------------------------------------------------------
$line = $wb.Worksheets.Item("line")

# Row 8 - Line6
$line.Range("A7").Copy($line.Range("A8"))
$line.Range("A8").Value2 = 6
$line.Range("B8").Value2 = "Line6"
$line.Range("C8").Value2 = "double"
$line.Range("D8").Value2 = 4
$line.Range("E8").Value2 = 5
$line.Range("F8").Value2 = 136
$line.Range("G8").Value2 = 0.031
$line.Range("H8").Value2 = 0.191
$line.Range("I8").Value2 = 18.8904
$line.Range("J8").Value2 = 0
$line.Range("K8").Value2 = 1.778
$line.Range("L8").Value2 = 1
$line.Range("M8").Value2 = 1
$line.Range("O8").Value2 = $true

# Row 9 - Line7
$line.Range("A7").Copy($line.Range("A9"))
$line.Range("A9").Value2 = 7
$line.Range("B9").Value2 = "Line7"
$line.Range("C9").Value2 = "double"
$line.Range("D9").Value2 = 4
$line.Range("E9").Value2 = 1
$line.Range("F9").Value2 = 154
$line.Range("G9").Value2 = 0.031
$line.Range("H9").Value2 = 0.191
$line.Range("I9").Value2 = 18.8904
$line.Range("J9").Value2 = 0
$line.Range("K9").Value2 = 1.778
$line.Range("L9").Value2 = 1
$line.Range("M9").Value2 = 1
$line.Range("O9").Value2 = $true

# ---------------------------------------------------------------------------
# Sheet "trafo": upgrade the 100MVA transformers (rows 2, 4, 5) to 200MVA
# (the 200MVA transformer already in row 3 is used as the reference values)
# ---------------------------------------------------------------------------
$trafo = $wb.Worksheets.Item("trafo")

$trafoRows = @(2, 4, 5)
foreach ($r in $trafoRows) {
    $trafo.Range("C$r").Value2 = "200MVA"
    $trafo.Range("F$r").Value2 = 200
    $trafo.Range("I$r").Value2 = 12.2
    $trafo.Range("J$r").Value2 = 0.26
    $trafo.Range("K$r").Value2 = 65
    $trafo.Range("L$r").Value2 = 0.06
    $trafo.Range("M$r").Value2 = 0
}
